$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.2
$ws.Range("G2").Value = 1.22
$ws.Range("H2").Value = 22
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 7.2
$ws.Range("K2").Value = 7.8
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 6.6
$ws.Range("O2").Value = 1.16
$ws.Range("P2").Value = 2.32
$ws.Range("Q2").Value = 1.71
$ws.Range("R2").Value = 1.42
$ws.Range("S2").Value = 3.2
$ws.Range("T2").Value = 2.06
$ws.Range("U2").Value = 1.76
$ws.Range("V2").Value = 1.04
$ws.Range("W2").Value = 5.5
$ws.Range("Z2").Value = 1000
$ws.Range("AB2").Value = 7.6
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 55
$ws.Range("AE2").Value = 280
$ws.Range("AF2").Value = 5.9
$ws.Range("AG2").Value = 9
$ws.Range("AH2").Value = 36
$ws.Range("AI2").Value = 160
$ws.Range("AJ2").Value = 8.4
$ws.Range("AK2").Value = 14.5
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 280
$ws.Range("AN2").Value = 7
